# Updates the cryptos price/volume table (and one B/C row swap) to match
# the refreshed GitHub Actions data pull. Cells hold plain text (not numbers)
# in the source workbook (e.g. "89.828.26", "1.00", "  -0.92%  "), so each
# target cell is forced to Text format ("@") before the value is written --
# otherwise Excel would auto-coerce number-looking strings (like "0.999" or
# "0.680") into real numbers and silently drop meaningful trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.828.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.92%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.95%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.14%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.43"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.48%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.71%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.361"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.59%  "
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.07%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.071.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.97%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.732"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.31%  "
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.66%  "
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.80%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.01%  "
# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.05%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.693.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.71%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.653.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.08%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.076.48"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.06%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.18%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.31%  "
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.92%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.79%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "434.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.34%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.62"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.07%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.33%  "
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.34%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.227.93"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.64%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +20.15%  "
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.21%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +30.31%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.14%  "
# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.165"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.24%  "
# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.966"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.36%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.76"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.94%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +20.24%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.05"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.33%  "
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.37%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "483.11"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.55%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -9.49%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.63%  "
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.48%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.14"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.25%  "
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.00%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.93"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.19%  "
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.46%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.680"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.76%  "
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.48%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.45%  "
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.23%  "
